$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid (G) and Absent (H) both become 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-18: Absent (H) becomes 1 for each day's row
for ($row = 4; $row -le 18; $row++) {
    $ws.Cells.Item($row, 8).Value = 1
}
